# Auto-generated PowerShell COM-interop script to apply cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.973.26"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.639.89"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D5").Value = "'213.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'23.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").Value = "1.872.98"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.642.52"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "'4.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("D16").Value = "'65.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "27.969.53"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "'232.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'7.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'10.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'4.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'151.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").Value = "'15.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").Value = "1.406.38"
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'0.882"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D39").Value = "'0.557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  -5.06%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +7.09%  "
$ws.Range("D44").Value = "'66.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "'5.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "1.781.62"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "'88.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'7.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
